# Mise a jour de l'application
# Adds two new training-session date columns (DF, DG = 2025-11-03 and
# 2025-11-04, serials 46048/46049) to the attendance sheet, fills in the
# attendance letters for existing players, fixes three previously-wrong
# "B" (Blessure) marks for row 2 that should have been "P" (Present),
# and appends a brand-new player row ("Mehdi Boussaid") with his own
# running-total formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row: two new date columns
# ---------------------------------------------------------------------
$ws.Range("DF1").Value = 46048
$ws.Range("DG1").Value = 46049
$ws.Range("DE1").Copy()
$ws.Range("DF1:DG1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Correct three stray "B" marks in row 2 that should read "P"
# ---------------------------------------------------------------------
$ws.Range("DC2").Value = "P"
$ws.Range("DD2").Value = "P"
$ws.Range("DE2").Value = "P"

# ---------------------------------------------------------------------
# 3) Attendance letters for the two new dates, existing players
#    (row number -> [DF value, DG value])
# ---------------------------------------------------------------------
$attendance = @{
    2  = @("P","P")
    3  = @("P","P")
    4  = @("P","P")
    5  = @("P","P")
    6  = @("A","P")
    7  = @("P","P")
    8  = @("P","P")
    9  = @("P","P")
    10 = @("P","P")
    11 = @("P","P")
    13 = @("RH","RH")
    14 = @("P","P")
    15 = @("P","P")
    16 = @("P","P")
    18 = @("B","RH")
    19 = @("P","P")
    20 = @("P","P")
    22 = @("P","P")
    24 = @("P","P")
    26 = @("A","A")
    27 = @("P","P")
    28 = @("P","P")
    29 = @("P","P")
    30 = @("P","P")
}

foreach ($r in $attendance.Keys) {
    $vals = $attendance[$r]
    $ws.Cells.Item($r, 110).Value = $vals[0]
    $ws.Cells.Item($r, 111).Value = $vals[1]
    $ws.Cells.Item($r, 109).Copy()
    $pasteRange = $ws.Range($ws.Cells.Item($r, 110), $ws.Cells.Item($r, 111))
    $pasteRange.PasteSpecial(-4122)
}

# Rows 17 and 25 stop earlier (players no longer tracked): the two new
# columns still need blank, styled cells (same as their existing empty
# DB:DE cells), with no value.
foreach ($r in @(17, 25)) {
    $ws.Cells.Item($r, 109).Copy()
    $pasteRange = $ws.Range($ws.Cells.Item($r, 110), $ws.Cells.Item($r, 111))
    $pasteRange.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 4) New player row 31 ("Mehdi Boussaid") with his own summary formulas
# ---------------------------------------------------------------------
$ws.Range("A31").Value = "Mehdi Boussaid"
$ws.Range("B31").Formula = "=COUNTA(K31:VR31)"
$ws.Range("C31").Formula = "=COUNTIF(K31:VR31,""P"")"
$ws.Range("D31").Formula = "=COUNTIF(K31:VR31,""REP"")"
$ws.Range("E31").Formula = "=COUNTIF(K31:VS31,""A"")"
$ws.Range("F31").Formula = "=COUNTIF(K31:VT31,""B"")"
$ws.Range("G31").Formula = "=COUNTIF(K31:VU31,""M"")"
$ws.Range("H31").Formula = "=COUNTIF(K31:VV31,""R"")"
$ws.Range("I31").Formula = "=COUNTIF(K31:VW31,""RH"")"
$ws.Range("J31").Formula = "=COUNTIF(K31:VW31,""S"")"
$ws.Range("DF31").Value = "P"
$ws.Range("DG31").Value = "P"

$ws.Range("A30:J30").Copy()
$ws.Range("A31:J31").PasteSpecial(-4122)
$ws.Range("DE30").Copy()
$ws.Range("DF31:DG31").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5) Conditional formatting range grows by one row
# ---------------------------------------------------------------------
$fcs = $ws.Range("A21:A30").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A21:A31"))
}

# ---------------------------------------------------------------------
# 6) Selection / view state
# ---------------------------------------------------------------------
$ws.Range("DI29").Select()

$excel.CutCopyMode = $false
